$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.348.62'
$ws.Range('E2').Value = '  +4.25%  '

$ws.Range('D3').Value = '1.548.27'
$ws.Range('E3').Value = '  +5.18%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('E5').Value = '  +0.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '282.65'
$ws.Range('E6').Value = '  +2.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3626'
$ws.Range('E7').Value = '  -0.70%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3202'
$ws.Range('E8').Value = '  +4.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.93'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').Value = '  +6.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06930'
$ws.Range('E11').Value = '  +4.79%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.13%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.730'
$ws.Range('E13').Value = '  +5.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.96'
$ws.Range('E14').Value = '  +3.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.437'
$ws.Range('E15').Value = '  +4.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001055'
$ws.Range('E16').Value = '  +2.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9685'
$ws.Range('E17').Value = '  -0.55%  '

$ws.Range('D18').Value = '1.544.24'
$ws.Range('E18').Value = '  +4.81%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06150'
$ws.Range('E19').Value = '  +4.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.24'
$ws.Range('E20').Value = '  +6.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.756'
$ws.Range('E21').Value = '  +5.64%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.31'
$ws.Range('E22').Value = '  +6.33%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.43'
$ws.Range('E23').Value = '  +4.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.321'
$ws.Range('E24').Value = '  +3.21%  '

$ws.Range('D25').Value = '21.360.07'
$ws.Range('E25').Value = '  +4.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.93'
$ws.Range('E26').Value = '  +4.34%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.286'
$ws.Range('E27').Value = '  +7.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.81'
$ws.Range('E28').Value = '  +3.43%  '

$ws.Range('D29').Value = '1.715.58'
$ws.Range('E29').Value = '  +5.39%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.67'
$ws.Range('E30').Value = '  +4.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.055'
$ws.Range('E31').Value = '  +4.51%  '

$ws.Range('E32').Value = '  +8.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.283'
$ws.Range('E33').Value = '  +6.62%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08056'
$ws.Range('E34').Value = '  +2.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.536'
$ws.Range('E35').Value = '  +0.47%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.002'
$ws.Range('E36').Value = '  +5.26%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.215'
$ws.Range('E37').Value = '  -1.99%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05874'
$ws.Range('E38').Value = '  +2.52%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2006'
$ws.Range('E39').Value = '  +6.76%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.019'
$ws.Range('E40').Value = '  +4.92%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02130'
$ws.Range('E41').Value = '  +4.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.82'
$ws.Range('E42').Value = '  +3.81%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9690'
$ws.Range('E43').Value = '  -0.07%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5531'
$ws.Range('E44').Value = '  +4.65%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.55'
$ws.Range('E45').Value = '  +4.19%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.585'
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5514'
$ws.Range('E47').Value = '  +6.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.31'
$ws.Range('E48').Value = '  +4.53%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.886'
$ws.Range('E49').Value = '  +6.69%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06624'
$ws.Range('E50').Value = '  +2.76%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.21'
$ws.Range('E51').Value = '  +5.09%  '

# Reset styles back to Normal for cells that required temporary text format
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
